$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Appearance" text for the Metal Gate (row 5) and Metal Railing (row 6) rows
$ws.Range("C5").Value = "Worn metal bars, making a net of pipes."
$ws.Range("C6").Value = "Old metal railings, modular."

# Update the active cell selection to reflect where the author last edited
$ws.Range("C6").Select()
